# Insert a new price-record row for "Vega Modelo de Temuco - Perejil" just
# above the existing row 439, pushing the remaining weekly records (formerly
# rows 439-460) down by one row to 440-461. This mirrors the author's commit
# "Fruta / hortaliza, semanal" which adds this week's new observation to the
# top of the dated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 439..460 down to 440..461, inheriting row 439's
# formatting (date number format, etc.) for the new row.
$ws.Rows.Item(439).Insert()

# Populate the newly inserted row 439 with the new weekly record.
$ws.Cells.Item(439, 1).Value = 10
$ws.Cells.Item(439, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(439, 3).Value = "La Araucanía"
$ws.Cells.Item(439, 4).Value = "2023-03-23"
$ws.Cells.Item(439, 5).Value = 9
$ws.Cells.Item(439, 6).Value = 100112044
$ws.Cells.Item(439, 7).Value = "Perejil"
$ws.Cells.Item(439, 8).Value = "Sin especificar"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 65
$ws.Cells.Item(439, 11).Value = 4000
$ws.Cells.Item(439, 12).Value = 4000
$ws.Cells.Item(439, 13).Value = 4000
$ws.Cells.Item(439, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(439, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(439, 16).Value = 1333
$ws.Cells.Item(439, 17).Value = 3
$ws.Cells.Item(439, 18).Value = "Hortaliza"
